$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.598.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.237.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "271.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.37%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.615"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0920"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.569.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.231.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.796"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.599.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("E19").Value = "  -1.05%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.16%  "

$ws.Range("E22").Value = "  -0.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.28%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.24%  "

$ws.Range("E28").Value = "  +5.08%  "

$ws.Range("E29").Value = "  +2.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("E32").Value = "  +1.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.06%  "

$ws.Range("E35").Value = "  +0.24%  "

$ws.Range("E36").Value = "  -4.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0350"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.32%  "

$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.213"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.05%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0990"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.51%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.14%  "

$ws.Range("E48").Value = "  +2.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.442"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.54%  "
